$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text-formatted, matching the source data
# (values like "1.000" / "0.000006753" / "8.730" must not be coerced to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.537.26'
$ws.Range("E2").Value = '  +2.06%  '
$ws.Range("D3").Value = '1.666.97'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '238.26'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D7").Value = '0.4795'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").Value = '0.2633'
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").Value = '0.06169'
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("D10").Value = '0.07096'
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").Value = '1.663.33'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '14.83'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '0.5903'
$ws.Range("E13").Value = '  -5.02%  '
$ws.Range("D14").Value = '4.372'
$ws.Range("E14").Value = '  -4.91%  '
$ws.Range("D15").Value = '75.11'
$ws.Range("E15").Value = '  +3.05%  '
$ws.Range("D16").Value = '0.9994'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '25.543.40'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("D19").Value = '0.000006753'
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '1.874.13'
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").Value = '4.422'
$ws.Range("E22").Value = '  -2.90%  '
$ws.Range("D23").Value = '8.730'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '5.282'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").Value = '135.84'
$ws.Range("E25").Value = '  +2.70%  '
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").Value = '1.388'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '105.14'
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("E29").Value = '  +2.64%  '
$ws.Range("D30").Value = '3.976'
$ws.Range("E30").Value = '  +5.72%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.07743'
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.643'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").Value = '0.9986'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '0.04236'
$ws.Range("E34").Value = '  -7.87%  '
$ws.Range("D35").Value = '2.598'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '0.6121'
$ws.Range("E36").Value = '  +6.03%  '
$ws.Range("D37").Value = '0.9504'
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("D38").Value = '2.594'
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("D39").Value = '0.8644'
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("D40").Value = '0.9993'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '1.852'
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  -6.04%  '
$ws.Range("D43").Value = '97.19'
$ws.Range("E43").Value = '  -0.94%  '
$ws.Range("D44").Value = '0.3767'
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("D45").Value = '4.855'
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '0.1122'
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").Value = '6.210'
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("D48").Value = '0.05266'
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("D49").Value = '29.77'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.384'
$ws.Range("E50").Value = '  +2.75%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.01%  '
